$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7992891073226929
$ws.Range("B1").Value = 2.767135381698608
$ws.Range("C1").Value = 7.643196105957031
$ws.Range("D1").Value = 2.244378805160522
$ws.Range("E1").Value = 1.484423398971558
